$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 277.93332
$ws.Range("I6").Value = 103.6
$ws.Range("J6").Value = 365.1
$ws.Range("K6").Value = 310.8
$ws.Range("L6").Value = 1095.3
$ws.Range("M6").Value = -198.8
$ws.Range("N6").Value = -1319.3

$ws.Range("H9").Value = 221.30435
$ws.Range("I9").Value = 142.25
$ws.Range("J9").Value = 307.54544
$ws.Range("K9").Value = 142.25
$ws.Range("L9").Value = 307.54544
$ws.Range("M9").Value = 26.75
$ws.Range("N9").Value = -645.54544

$ws.Range("H17").Value = 702
$ws.Range("J17").Value = 675.1818
$ws.Range("L17").Value = 2025.5454
$ws.Range("N17").Value = -2361.5454

$ws.Range("H113").Value = 8732.637000000001
$ws.Range("I113").Value = 7105.9
$ws.Range("J113").Value = 25000
$ws.Range("K113").Value = 7105.9
$ws.Range("L113").Value = 25000
$ws.Range("M113").Value = -3851.9
$ws.Range("N113").Value = -31508

$ws.Range("H116").Value = 6399.478
$ws.Range("I116").Value = 4085.5
$ws.Range("J116").Value = 7633.6
$ws.Range("K116").Value = 4085.5
$ws.Range("L116").Value = 7633.6
$ws.Range("M116").Value = -643.5
$ws.Range("N116").Value = -14517.6

$ws.Range("H132").Value = 19510.516
$ws.Range("I132").Value = 2013.12
$ws.Range("K132").Value = 6039.36
$ws.Range("M132").Value = -3509.36

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").ClearContents()
$ws.Range("N104").Value = 0

$ws.Range("H122").Value = 4396.4834
$ws.Range("I122").Value = 3730.9375
$ws.Range("J122").Value = 5157.107
$ws.Range("K122").Value = 11192.8125
$ws.Range("L122").Value = 15471.321
$ws.Range("M122").Value = -8742.8125
$ws.Range("N122").Value = -20371.321

$ws.Range("H131").Value = 34745.875
$ws.Range("J131").Value = 34745.875
$ws.Range("L131").Value = 34745.875
$ws.Range("N131").Value = -44825.875

$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 36139.137
$ws.Range("I86").Value = 41486.72
$ws.Range("J86").Value = 2716.75
$ws.Range("K86").Value = 41486.72
$ws.Range("L86").Value = 2716.75
$ws.Range("M86").Value = -40363.72
$ws.Range("N86").Value = -4962.75

$ws.Range("H89").Value = 36139.137
$ws.Range("I89").Value = 41486.72
$ws.Range("J89").Value = 2716.75
$ws.Range("K89").Value = 207433.6
$ws.Range("L89").Value = 13583.75
$ws.Range("M89").Value = -201817.6
$ws.Range("N89").Value = -24815.75

$ws.Range("H134").Value = 3123.353
$ws.Range("I134").Value = 3074.8125
$ws.Range("J134").Value = 3900
$ws.Range("K134").Value = 9224.4375
$ws.Range("L134").Value = 11700
$ws.Range("M134").Value = -6689.4375
$ws.Range("N134").Value = -16770

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 0

$ws.Range("H105").Value = 1870.2222
$ws.Range("I105").Value = 1854
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 1854
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -107
$ws.Range("N105").Value = -5494

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6332.6665
$ws.Range("I56").Value = 6332.6665
$ws.Range("K56").Value = 6332.6665
$ws.Range("M56").Value = -5802.6665

$ws.Range("H68").Value = 2857.6
$ws.Range("I68").Value = 642.5
$ws.Range("K68").Value = 1927.5
$ws.Range("M68").Value = -1116.5

$ws.Range("H71").Value = 2857.6
$ws.Range("I71").Value = 642.5
$ws.Range("K71").Value = 5782.5
$ws.Range("M71").Value = -1726.5

$ws.Range("H102").Value = 4000
$ws.Range("J102").Value = 4000
$ws.Range("L102").Value = 12000
$ws.Range("N102").Value = -16868

$ws.Range("H105").Value = 8000
$ws.Range("J105").Value = 8000
$ws.Range("L105").Value = 24000
$ws.Range("N105").Value = -29242

$ws.Range("H107").Value = 2166.3547
$ws.Range("J107").Value = 2080.7827
$ws.Range("L107").Value = 6242.348100000001
$ws.Range("N107").Value = -10082.3481

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 90492.62
$ws.Range("I70").Value = 143249.88
$ws.Range("J70").Value = 6081
$ws.Range("K70").Value = 143249.88
$ws.Range("L70").Value = 6081
$ws.Range("M70").Value = -142979.88
$ws.Range("N70").Value = -6621

$ws.Range("H73").Value = 90492.62
$ws.Range("I73").Value = 143249.88
$ws.Range("J73").Value = 6081
$ws.Range("K73").Value = 143249.88
$ws.Range("L73").Value = 6081
$ws.Range("M73").Value = -142313.88
$ws.Range("N73").Value = -7953

$ws.Range("H102").Value = 3024.3076
$ws.Range("I102").Value = 2463.8333
$ws.Range("K102").Value = 2463.8333
$ws.Range("M102").Value = -841.8332999999998

$ws.Range("H122").Value = 3586.9688
$ws.Range("I122").Value = 2964
$ws.Range("J122").Value = 4625.25
$ws.Range("K122").Value = 8892
$ws.Range("L122").Value = 13875.75
$ws.Range("M122").Value = -6442
$ws.Range("N122").Value = -18775.75

$ws.Range("H132").Value = 5383.569
$ws.Range("I132").Value = 4991.061
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 14973.183
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -12443.183
$ws.Range("N132").Value = -50060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1066.6666
$ws.Range("I22").Value = 1066.6666
$ws.Range("K22").Value = 1066.6666
$ws.Range("M22").Value = -771.6666

$ws.Range("H27").Value = 1066.6666
$ws.Range("I27").Value = 1066.6666
$ws.Range("K27").Value = 1066.6666
$ws.Range("M27").Value = -959.6666

$ws.Range("H40").Value = 3643.7368
$ws.Range("I40").Value = 3324.7334
$ws.Range("J40").Value = 4840
$ws.Range("K40").Value = 3324.7334
$ws.Range("L40").Value = 4840
$ws.Range("M40").Value = -3188.7334
$ws.Range("N40").Value = -5112

$ws.Range("H122").Value = 4065.7727
$ws.Range("I122").Value = 3349.8235
$ws.Range("J122").Value = 6500
$ws.Range("K122").Value = 10049.4705
$ws.Range("L122").Value = 19500
$ws.Range("M122").Value = -7599.470499999999
$ws.Range("N122").Value = -24400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 19909.092
$ws.Range("J57").Value = 19909.092
$ws.Range("L57").Value = 19909.092
$ws.Range("N57").Value = -21417.092

$ws.Range("H109").Value = 14000
$ws.Range("J109").Value = 14000
$ws.Range("L109").Value = 14000
$ws.Range("N109").Value = -16774

$ws.Range("H113").Value = 1096.6666
$ws.Range("I113").Value = 539.75
$ws.Range("J113").Value = 2210.5
$ws.Range("K113").Value = 1619.25
$ws.Range("L113").Value = 6631.5
$ws.Range("M113").Value = 550.75
$ws.Range("N113").Value = -10971.5

$ws.Range("H122").Value = 1239.0555
$ws.Range("I122").Value = 1138.6923
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 3416.0769
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -966.0769
$ws.Range("N122").Value = -9400

$ws.Range("H132").Value = 3217.125
$ws.Range("I132").Value = 3277.2307
$ws.Range("J132").Value = 2956.6667
$ws.Range("K132").Value = 9831.6921
$ws.Range("L132").Value = 8870.000100000001
$ws.Range("M132").Value = -7301.6921
$ws.Range("N132").Value = -13930.0001
